$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Opportunity IDs (column A) are stored as text in the source data, not
# numbers -- force text formatting before writing the numeric-looking
# strings so they don't get auto-converted to numbers.
$ws.Range("A2:A4").NumberFormat = "@"

# Update row 2
$ws.Range("A2").Value = "1328650"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328650"
$ws.Range("C2").Value = "Marketing Intern"
$ws.Range("D2").Value = "Novi Sad, Serbia"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "2 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "DataDrill"

# Update row 3
$ws.Range("A3").Value = "1328649"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328649"
$ws.Range("C3").Value = "Sales Development Representative"
$ws.Range("D3").Value = "Novi Sad, Serbia"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "DataDrill"

# Add new row 4
$ws.Range("A4").Value = "1328516"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1328516"
$ws.Range("C4").Value = "Account Manager"
$ws.Range("D4").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "Tech-Pix"

# Column width adjustments.
# Excel's ColumnWidth (COM) is offset from the stored OOXML <col width>
# by the default font's padding (~0.8333 chars for Calibri 11), so back
# that padding out to land on the exact stored widths from the diff.
$padding = 0.8333333333333333
$ws.Columns.Item(3).ColumnWidth = 35 - $padding
$ws.Columns.Item(4).ColumnWidth = 34 - $padding
$ws.Columns.Item(6).ColumnWidth = 15 - $padding
$ws.Columns.Item(7).ColumnWidth = 15 - $padding
$ws.Columns.Item(8).ColumnWidth = 15 - $padding
